$wb = $excel.ActiveWorkbook
$wsMeta = $wb.Worksheets.Item("Metadata")
$wsConcepts = $wb.Worksheets.Item("Concepts")

# --- Metadata sheet updates ---

# Date (plain text, not numeric-looking -> safe to set directly).
$wsMeta.Cells.Item(8, 2).Value2 = "2025-05-21T20:08:08+00:00"

# Count: "8" -> "10". "10" looks numeric, and setting it directly (via Value /
# Value2 / even a leading apostrophe or NumberFormat="@") makes the engine
# either store it as a real number (losing the shared-string "t=s" marker) or
# bump the cell onto a brand-new style slot (quotePrefix / explicit numFmt).
# Route the text through a scratch cell + Copy/PasteSpecial(values-only) so the
# destination keeps its original style (s="2") while still getting a text
# ("t=s") cell, exactly like the original "8" was stored.
$scratch = $wsMeta.Cells.Item(200, 200)
$scratch.Value = "'10"
$scratch.Copy()
$wsMeta.Cells.Item(22, 2).PasteSpecial(-4163)   # xlPasteValues
$scratch.Clear()

# --- Concepts sheet: insert two new rows while preserving existing data/styles ---
# Shift existing rows down using Range.Copy(destination), which carries both
# values (as proper shared strings) and styles, processed bottom-up so we
# never read a row after it has been overwritten.
$wsConcepts.Range("A9:D9").Copy($wsConcepts.Range("A11:D11"))   # Annotated-SNV row
$wsConcepts.Range("A8:D8").Copy($wsConcepts.Range("A10:D10"))   # IGV row
$wsConcepts.Range("A7:D7").Copy($wsConcepts.Range("A9:D9"))     # Germline-structural-variant row
$wsConcepts.Range("A6:D6").Copy($wsConcepts.Range("A8:D8"))     # Sequencing-data-index row
$wsConcepts.Range("A5:D5").Copy($wsConcepts.Range("A6:D6"))     # Sequencing-data-supplement row
$wsConcepts.Range("A4:D4").Copy($wsConcepts.Range("A5:D5"))     # Germline-CNV row
$wsConcepts.Range("A3:D3").Copy($wsConcepts.Range("A4:D4"))     # SNV row
$wsConcepts.Range("A2:D2").Copy($wsConcepts.Range("A3:D3"))     # Aligned-reads row

# Row 2 (A2/D2 already correctly "1" / blank from the original row) becomes the
# new "Unaligned-Reads" concept.
$wsConcepts.Cells.Item(2, 2).Value2 = "Unaligned-Reads"
$wsConcepts.Cells.Item(2, 3).Value2 = "Unaligned Reads"

# Row 7 (A7/D7 still correctly "1" / blank, never overwritten above) becomes the
# new "Metrics" concept.
$wsConcepts.Cells.Item(7, 2).Value2 = "Metrics"
$wsConcepts.Cells.Item(7, 3).Value2 = "Metrics"
